# Update "countries & provincias Spain" data (covid-19 country stats table).
# 1) Refresh the "last updated" timestamp in the title cell (A1).
# 2) Update the 7 numeric stat columns (B:H) for the countries whose figures
#    changed in this data refresh (by current row, before re-sorting).
# 3) Re-sort the whole data range (A4:H215) descending by column B
#    ("Casos totales"), exactly like the published sheet does after every
#    refresh, so rows whose totals grew (Pakistan, Peru, Monaco, Aruba, ...)
#    bubble back up above their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Title / timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 20:22"

# --- 2) Per-country numeric updates (row numbers as they exist right now,
#        i.e. before the re-sort below moves anything) ---------------------
$updates = @{
    4   = @(356414, 19741, 19247, 326677, 8857, 874, 10490)   # Estados Unidos
    7   = @(101558, 1435, 28700, 71196, 3936, 78, 1662)       # Alemania
    17  = @(12286, 235, 3463, 8603, 250, 16, 220)             # Austria
    25  = @(5763, 76, 32, 5655, 83, 5, 76)                    # Noruega
    28  = @(4778, 489, 375, 4267, 0, 18, 136)                 # India
    35  = @(3766, 609, 259, 3455, 17, 5, 52)                  # Pakistan
    41  = @(2561, 280, 997, 1472, 89, 9, 92)                  # Peru
    67  = @(843, 32, 8, 820, 11, 2, 15)                       # Lituania
    130 = @(77, 4, 4, 72, 4, 0, 1)                            # Monaco
    133 = @(69, 7, 5, 60, 4, 1, 4)                            # El Salvador
    135 = @(71, 7, 2, 69, 0, 0, 0)                            # Aruba
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# --- 3) Re-sort the data rows (A4:H215) by "Casos totales" (col B), desc --
$dataRange = $ws.Range("A4:H215")
$sortKey = $ws.Range("B4:B215")
$dataRange.Sort($sortKey, 2)

# Monaco's updated total (77) ties with Liechtenstein's unchanged total (77).
# The source refresh places Monaco immediately above Liechtenstein on a tie;
# our stable sort instead keeps Liechtenstein's earlier position, so swap the
# two rows back into the published order.
for ($r = 4; $r -le 214; $r++) {
    $nameHere = [string]$ws.Cells.Item($r, 1).Value()
    $nameNext = [string]$ws.Cells.Item($r + 1, 1).Value()
    if ($nameHere -eq "Liechtenstein" -and $nameNext -eq "Monaco") {
        for ($c = 1; $c -le 8; $c++) {
            $v1 = $ws.Cells.Item($r, $c).Value()
            $v2 = $ws.Cells.Item($r + 1, $c).Value()
            $ws.Cells.Item($r, $c).Value = $v2
            $ws.Cells.Item($r + 1, $c).Value = $v1
        }
    }
}
